$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 18 new rows before row 2 to make room for the 6/08 moisture data.
# This shifts the existing rows (formerly 2-127) down to rows 20-145.
$ws.Rows("2:19").Insert()

# Fill in the new rows (2-19) with the 6/08 (soil_type x pot) data.
# Column D (sd) is intentionally left blank for these new rows.
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = 1
$ws.Range("C2").Value = 47.5
$ws.Range("E2").Value = 44779
$ws.Range("A3").Value = 1
$ws.Range("B3").Value = 2
$ws.Range("C3").Value = 50.2
$ws.Range("E3").Value = 44779
$ws.Range("A4").Value = 1
$ws.Range("B4").Value = 3
$ws.Range("C4").Value = 47.3
$ws.Range("E4").Value = 44779
$ws.Range("A5").Value = 1
$ws.Range("B5").Value = 4
$ws.Range("C5").Value = 51.5
$ws.Range("E5").Value = 44779
$ws.Range("A6").Value = 1
$ws.Range("B6").Value = 5
$ws.Range("C6").Value = 51.3
$ws.Range("E6").Value = 44779
$ws.Range("A7").Value = 1
$ws.Range("B7").Value = 6
$ws.Range("C7").Value = 53.2
$ws.Range("E7").Value = 44779
$ws.Range("A8").Value = 2
$ws.Range("B8").Value = 1
$ws.Range("C8").Value = 49
$ws.Range("E8").Value = 44779
$ws.Range("A9").Value = 2
$ws.Range("B9").Value = 2
$ws.Range("C9").Value = 55.3
$ws.Range("E9").Value = 44779
$ws.Range("A10").Value = 2
$ws.Range("B10").Value = 3
$ws.Range("C10").Value = 54.9
$ws.Range("E10").Value = 44779
$ws.Range("A11").Value = 2
$ws.Range("B11").Value = 4
$ws.Range("C11").Value = 50.8
$ws.Range("E11").Value = 44779
$ws.Range("A12").Value = 2
$ws.Range("B12").Value = 5
$ws.Range("C12").Value = 56.1
$ws.Range("E12").Value = 44779
$ws.Range("A13").Value = 2
$ws.Range("B13").Value = 6
$ws.Range("C13").Value = 55.6
$ws.Range("E13").Value = 44779
$ws.Range("A14").Value = 3
$ws.Range("B14").Value = 1
$ws.Range("C14").Value = 39.799999999999997
$ws.Range("E14").Value = 44779
$ws.Range("A15").Value = 3
$ws.Range("B15").Value = 2
$ws.Range("C15").Value = 45.1
$ws.Range("E15").Value = 44779
$ws.Range("A16").Value = 3
$ws.Range("B16").Value = 3
$ws.Range("C16").Value = 49.3
$ws.Range("E16").Value = 44779
$ws.Range("A17").Value = 3
$ws.Range("B17").Value = 4
$ws.Range("C17").Value = 38.299999999999997
$ws.Range("E17").Value = 44779
$ws.Range("A18").Value = 3
$ws.Range("B18").Value = 5
$ws.Range("C18").Value = 40.4
$ws.Range("E18").Value = 44779
$ws.Range("A19").Value = 3
$ws.Range("B19").Value = 6
$ws.Range("C19").Value = 50.6
$ws.Range("E19").Value = 44779

# Fix up the date column's style: Insert() copied the header row's
# formatting (bold-ish font, no number format) onto the new E cells.
# Reset to the plain "Normal" style and re-apply the date number format
# used by the rest of the date column (style index 2 / numFmtId 16).
$ws.Range("E2:E19").Style = "Normal"
$ws.Range("E2:E19").NumberFormat = "d-mmm"

# Update the active selection to match the author's final cursor position.
[void]$ws.Range("F19").Select()
